# Fixed naive component forecaster bug - Presentation state 11.02.
# Recomputed the naive QoQ error series matched to ifoCAST (AR2_50_9, T45):
# values shift down the "staircase" by one row (row 2 loses its last
# column K2, rows 15-24 each pick up one new trailing column).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Q10) loses its last value (old K2) after the recompute.
$ws.Range("K2").ClearContents()

$ws.Range("B2").Value = 2.45278688403014
$ws.Range("C2").Value = 10.34861343995435
$ws.Range("D2").Value = -7.968604821155352
$ws.Range("E2").Value = 0.1773919319351871
$ws.Range("F2").Value = 1.968570813843447
$ws.Range("G2").Value = -1.20814069825023
$ws.Range("H2").Value = -1.502781793147165
$ws.Range("I2").Value = 1.011999790233559
$ws.Range("J2").Value = 0.08244241780622588
$ws.Range("B3").Value = 1.537235586080065
$ws.Range("C3").Value = 9.793174525304376
$ws.Range("D3").Value = -8.287524531353117
$ws.Range("E3").Value = -0.0318104525324077
$ws.Range("F3").Value = 1.81535338203464
$ws.Range("G3").Value = -1.333595567856255
$ws.Range("H3").Value = -1.61435418494062
$ws.Range("I3").Value = 0.9073520362996292
$ws.Range("J3").Value = -0.01874886075804749
$ws.Range("K3").Value = 0.218064591657356
$ws.Range("B4").Value = 2.627944642717427
$ws.Range("C4").Value = -14.06217932890006
$ws.Range("D4").Value = -4.673868413464427
$ws.Range("E4").Value = -1.935174446597504
$ws.Range("F4").Value = -4.389680736670638
$ws.Range("G4").Value = -4.134058956105767
$ws.Range("H4").Value = -1.20012607464721
$ws.Range("I4").Value = -1.810454602675451
$ws.Range("J4").Value = -1.332263916244841
$ws.Range("K4").Value = -1.454243286040335
$ws.Range("B5").Value = -12.38199818233083
$ws.Range("C5").Value = 2.742176349044278
$ws.Range("D5").Value = 0.312881869627589
$ws.Range("E5").Value = -0.408319721481464
$ws.Range("F5").Value = -2.155877054195791
$ws.Range("G5").Value = 1.211880171423111
$ws.Range("H5").Value = -0.2204442241744874
$ws.Range("I5").Value = 0.3103086300649579
$ws.Range("J5").Value = -0.1709206463330128
$ws.Range("K5").Value = 0.4422206679346297
$ws.Range("B6").Value = -0.9850301776737167
$ws.Range("C6").Value = 2.175752524391212
$ws.Range("D6").Value = -1.169228645652444
$ws.Range("E6").Value = -1.605646058573539
$ws.Range("F6").Value = 0.9646269304998711
$ws.Range("G6").Value = 0.04013991644653381
$ws.Range("H6").Value = 0.2665276098352697
$ws.Range("I6").Value = -0.0392331488638116
$ws.Range("J6").Value = 0.4722505897783218
$ws.Range("K6").Value = 0.1845085656625323
$ws.Range("B7").Value = 2.583984978012982
$ws.Range("C7").Value = -0.9176056017660498
$ws.Range("D7").Value = -1.687552112396674
$ws.Range("E7").Value = 1.039718979881157
$ws.Range("F7").Value = 0.1460470769602576
$ws.Range("G7").Value = 0.3249442559939251
$ws.Range("H7").Value = 0.02720553496627992
$ws.Range("I7").Value = 0.5469325483891583
$ws.Range("J7").Value = 0.2548218575635111
$ws.Range("K7").Value = 0.2823415042266514
$ws.Range("B8").Value = -1.023104813614865
$ws.Range("C8").Value = -1.813657396362032
$ws.Range("D8").Value = 1.207967048514389
$ws.Range("E8").Value = 0.2067111695641399
$ws.Range("F8").Value = 0.336077095368558
$ws.Range("G8").Value = 0.09006828629321051
$ws.Range("H8").Value = 0.6042840761829015
$ws.Range("I8").Value = 0.2993270345678901
$ws.Range("J8").Value = 0.33341537525307
$ws.Range("K8").Value = 0.4282859431822899
$ws.Range("B9").Value = -2.418098016209771
$ws.Range("C9").Value = 1.143763671159695
$ws.Range("D9").Value = 0.4982584776213567
$ws.Range("E9").Value = 0.2983702192826244
$ws.Range("F9").Value = 0.08148426621974941
$ws.Range("G9").Value = 0.6880251537802698
$ws.Range("H9").Value = 0.3341541794959351
$ws.Range("I9").Value = 0.359979758305367
$ws.Range("J9").Value = 0.4737847839288888
$ws.Range("K9").Value = 0.6788178843816368
$ws.Range("B10").Value = 1.056122842320784
$ws.Range("C10").Value = 0.4458392755916881
$ws.Range("D10").Value = 0.35701661626763
$ws.Range("E10").Value = 0.08641341155685073
$ws.Range("F10").Value = 0.6769318707313746
$ws.Range("G10").Value = 0.3469208787210216
$ws.Range("H10").Value = 0.3687093173527261
$ws.Range("I10").Value = 0.4762867280777406
$ws.Range("J10").Value = 0.6850938218533179
$ws.Range("K10").Value = 0.1817482544735015
$ws.Range("B11").Value = 0.4682672552246967
$ws.Range("C11").Value = 0.3596201558300787
$ws.Range("D11").Value = 0.06434572590283549
$ws.Range("E11").Value = 0.6709498904270196
$ws.Range("F11").Value = 0.3415776285386105
$ws.Range("G11").Value = 0.3583407839823598
$ws.Range("H11").Value = 0.4678094059541168
$ws.Range("I11").Value = 0.6773215135465072
$ws.Range("J11").Value = 0.1731242809405329
$ws.Range("K11").Value = 0.4653169537672726
$ws.Range("B12").Value = 0.3810447272872988
$ws.Range("C12").Value = 0.1845235747020479
$ws.Range("D12").Value = 0.5817555460132365
$ws.Range("E12").Value = 0.3118506246900327
$ws.Range("F12").Value = 0.3669669040122699
$ws.Range("G12").Value = 0.4418692325286816
$ws.Range("H12").Value = 0.6545429505182796
$ws.Range("I12").Value = 0.1595062457919151
$ws.Range("J12").Value = 0.4467961565820968
$ws.Range("K12").Value = 0.1882892220727796
$ws.Range("B13").Value = 0.1422029423816584
$ws.Range("C13").Value = 0.5525510721575033
$ws.Range("D13").Value = 0.3088300248851115
$ws.Range("E13").Value = 0.3481056204007895
$ws.Range("F13").Value = 0.4215381014634773
$ws.Range("G13").Value = 0.6400127950840317
$ws.Range("H13").Value = 0.1428652002996471
$ws.Range("I13").Value = 0.4292123111533466
$ws.Range("J13").Value = 0.1718153385427368
$ws.Range("K13").Value = 0.5550469433309027
$ws.Range("B14").Value = 0.8927488127754134
$ws.Range("C14").Value = 0.3842010866690486
$ws.Range("D14").Value = 0.1570123340877904
$ws.Range("E14").Value = 0.4462320140502105
$ws.Range("F14").Value = 0.6493921986403277
$ws.Range("G14").Value = 0.08805440168825607
$ws.Range("H14").Value = 0.4103182202646768
$ws.Range("I14").Value = 0.1569513377276911
$ws.Range("J14").Value = 0.526335232701532
$ws.Range("K14").Value = 0.1198157454308724
$ws.Range("B15").Value = 0.8348959586592991
$ws.Range("C15").Value = 0.2041230363001488
$ws.Range("D15").Value = 0.2059487832594023
$ws.Range("E15").Value = 0.6853372797061905
$ws.Range("F15").Value = 0.08600364424866319
$ws.Range("G15").Value = 0.3337499808286466
$ws.Range("H15").Value = 0.1295519433524077
$ws.Range("I15").Value = 0.5000461825038066
$ws.Range("J15").Value = 0.0764184547175398
$ws.Range("K15").Value = 0.3733815035592049
$ws.Range("B16").Value = 0.5164486232236872
$ws.Range("C16").Value = 0.3392375229949266
$ws.Range("D16").Value = 0.506258857889999
$ws.Range("E16").Value = 0.1141978830192304
$ws.Range("F16").Value = 0.3697752920210401
$ws.Range("G16").Value = 0.09201037314819407
$ws.Range("H16").Value = 0.4945038434164454
$ws.Range("I16").Value = 0.08059378365628841
$ws.Range("J16").Value = 0.3621967898194339
$ws.Range("B17").Value = 0.5748280141027678
$ws.Range("C17").Value = 0.5909375987643086
$ws.Range("D17").Value = -0.03140248361448672
$ws.Range("E17").Value = 0.3797290348802828
$ws.Range("F17").Value = 0.107375542847739
$ws.Range("G17").Value = 0.4556555084590223
$ws.Range("H17").Value = 0.06474797462910251
$ws.Range("I17").Value = 0.3538625613519331
$ws.Range("B18").Value = 0.9019617852456914
$ws.Range("C18").Value = 0.08568629079670848
$ws.Range("D18").Value = 0.2172648934307159
$ws.Range("E18").Value = 0.1363344933129406
$ws.Range("F18").Value = 0.4913983809139259
$ws.Range("G18").Value = 0.03342032824547542
$ws.Range("H18").Value = 0.3506375116504074
$ws.Range("B19").Value = 0.3325114682008229
$ws.Range("C19").Value = 0.2349207609686054
$ws.Range("D19").Value = 0.0417122542019461
$ws.Range("E19").Value = 0.5236157691624059
$ws.Range("F19").Value = 0.0458084540731927
$ws.Range("G19").Value = 0.3297737081011285
$ws.Range("B20").Value = 0.4745971203848173
$ws.Range("C20").Value = 0.1267364976711596
$ws.Range("D20").Value = 0.4071311908043919
$ws.Range("E20").Value = 0.06091372572504519
$ws.Range("F20").Value = 0.3530127437271186
$ws.Range("B21").Value = 0.2915965747052469
$ws.Range("C21").Value = 0.4206684630523081
$ws.Range("D21").Value = -0.002059193264917797
$ws.Range("E21").Value = 0.3630509094086165
$ws.Range("B22").Value = 0.6745175049177161
$ws.Range("C22").Value = 0.09725885691711864
$ws.Range("D22").Value = 0.2514743671933744
$ws.Range("B23").Value = 0.1413185481565676
$ws.Range("C23").Value = 0.2728948720679752
$ws.Range("B24").Value = 0.5091380033804217
